# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 632
$ws1.Range("F3").Value = 208
$ws1.Range("F4").Value = 661
$ws1.Range("F6").Value = 316
$ws1.Range("F7").Value = 2811
$ws1.Range("F8").Value = 471
$ws1.Range("F9").Value = 7926
$ws1.Range("F13").Value = 378

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 632
$ws4.Range("F3").Value = 208
$ws4.Range("F4").Value = 661
$ws4.Range("F6").Value = 316
$ws4.Range("F9").Value = 2811
$ws4.Range("F10").Value = 471
$ws4.Range("F11").Value = 7926
$ws4.Range("F17").Value = 378
